$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.604.73"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.923.72"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.45"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4739"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06859"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "105.34"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.42"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").Value = "1.919.91"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07720"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.352"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6709"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "288.53"
$ws.Range("E16").Value = "  -6.40%  "
$ws.Range("D17").Value = "30.628.72"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007657"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.544"
$ws.Range("E21").Value = "  +4.03%  "
$ws.Range("D22").Value = "2.170.82"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.449"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.525"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.65"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.132"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.059"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05023"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7343"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.733"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.686"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "111.57"
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.047"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4443"
$ws.Range("E42").Value = "  +6.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8763"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.905"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.81"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.329"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.362"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1253"
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "48.05"
$ws.Range("E50").Value = "  +12.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.08"
$ws.Range("E51").Value = "  +0.54%  "
